$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the slightly-corrected timestamp value in A13
$ws.Range("A13").Value = 44326.78153125347

# Add the new row of data retrieved on 2021-05-11
$ws.Range("A14").Value = 44327.7819216314
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat

$ws.Range("B14").Value = 74185
$ws.Range("C14").Value = 62361
$ws.Range("D14").Value = 3242
$ws.Range("E14").Value = 2063
$ws.Range("F14").Value = 1458
$ws.Range("G14").Value = 19170
$ws.Range("H14").Value = 1328
$ws.Range("I14").Value = 848
$ws.Range("J14").Value = 215
